$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Shift old columns N:P (Late / heading / Outstanding) one slot to the right
# so a brand-new blank column can be inserted at N.
$ws.Range("N1:N14").Insert(-4161)

# Give the freshly inserted column the same width as column M and mark the
# data cells (rows 2-14) with the wrapped / vertically centred look used
# throughout the sheet (same formatting as the neighbouring data columns).
$ws.Columns("N").ColumnWidth = 10.7109375
$ws.Range("N2:N14").VerticalAlignment = -4108
$ws.Range("N2:N14").WrapText = $true

# --- Switch the active tab back to "Repayment schedule" (this also drops
# "Acc_Repayment1" as the active tab, leaving its own remembered selection
# untouched) and reposition the selection on it ---
$ws.Activate()
$ws.Range("J17").Select()
